# Auto-generated edit script for Asura_Profits (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
# Applies per-cell numeric updates (currentAveragePrice* / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 19
$ws_ALC.Cells.Item(19, 8).Value = 628.4167
$ws_ALC.Cells.Item(19, 10).Value = 643.875
$ws_ALC.Cells.Item(19, 12).Value = 643.875
$ws_ALC.Cells.Item(19, 14).Value = -993.875

# ALC row 70
$ws_ALC.Cells.Item(70, 8).Value = 7985319.5
$ws_ALC.Cells.Item(70, 9).Value = 33534186
$ws_ALC.Cells.Item(70, 10).Value = 1298.3125
$ws_ALC.Cells.Item(70, 11).Value = 100602558
$ws_ALC.Cells.Item(70, 12).Value = 3894.9375
$ws_ALC.Cells.Item(70, 13).Value = -100602288
$ws_ALC.Cells.Item(70, 14).Value = -4434.9375

# ALC row 73
$ws_ALC.Cells.Item(73, 8).Value = 7985319.5
$ws_ALC.Cells.Item(73, 9).Value = 33534186
$ws_ALC.Cells.Item(73, 10).Value = 1298.3125
$ws_ALC.Cells.Item(73, 11).Value = 100602558
$ws_ALC.Cells.Item(73, 12).Value = 3894.9375
$ws_ALC.Cells.Item(73, 13).Value = -100601622
$ws_ALC.Cells.Item(73, 14).Value = -5766.9375

# ALC row 128
$ws_ALC.Cells.Item(128, 8).Value = 23666.666
$ws_ALC.Cells.Item(128, 10).Value = 23666.666
$ws_ALC.Cells.Item(128, 12).Value = 23666.666
$ws_ALC.Cells.Item(128, 14).Value = -33626.666

# ARM row 44
$ws_ARM.Cells.Item(44, 8).Value = 20000
$ws_ARM.Cells.Item(44, 9).Value = 20000
$ws_ARM.Cells.Item(44, 11).Value = 20000
$ws_ARM.Cells.Item(44, 13).Value = -19512

# ARM row 55
$ws_ARM.Cells.Item(55, 8).Value = 22500

# ARM row 63
$ws_ARM.Cells.Item(63, 8).Value = 3572.1428
$ws_ARM.Cells.Item(63, 9).Value = 2601
$ws_ARM.Cells.Item(63, 11).Value = 2601
$ws_ARM.Cells.Item(63, 13).Value = -1915

# ARM row 66
$ws_ARM.Cells.Item(66, 8).Value = 3572.1428
$ws_ARM.Cells.Item(66, 9).Value = 2601
$ws_ARM.Cells.Item(66, 11).Value = 13005
$ws_ARM.Cells.Item(66, 13).Value = -9573

# ARM row 80
$ws_ARM.Cells.Item(80, 8).Value = 35633.332
$ws_ARM.Cells.Item(80, 10).Value = 35633.332
$ws_ARM.Cells.Item(80, 12).Value = 35633.332
$ws_ARM.Cells.Item(80, 14).Value = -37629.332

# ARM row 83
$ws_ARM.Cells.Item(83, 8).Value = 35633.332
$ws_ARM.Cells.Item(83, 10).Value = 35633.332
$ws_ARM.Cells.Item(83, 12).Value = 106899.996
$ws_ARM.Cells.Item(83, 14).Value = -116883.996

# BSM row 35
$ws_BSM.Cells.Item(35, 8).Value = 0
$ws_BSM.Cells.Item(35, 10).Value = 0
$ws_BSM.Cells.Item(35, 12).Value = 0
$ws_BSM.Cells.Item(35, 14).ClearContents()

# BSM row 82
$ws_BSM.Cells.Item(82, 8).Value = 62126.168
$ws_BSM.Cells.Item(82, 9).Value = 62126.168
$ws_BSM.Cells.Item(82, 11).Value = 62126.168
$ws_BSM.Cells.Item(82, 13).Value = -61743.168

# BSM row 85
$ws_BSM.Cells.Item(85, 8).Value = 62126.168
$ws_BSM.Cells.Item(85, 9).Value = 62126.168
$ws_BSM.Cells.Item(85, 11).Value = 62126.168
$ws_BSM.Cells.Item(85, 13).Value = -60800.168

# BSM row 94
$ws_BSM.Cells.Item(94, 8).Value = 1077.174
$ws_BSM.Cells.Item(94, 9).Value = 639.2
$ws_BSM.Cells.Item(94, 10).Value = 1898.375
$ws_BSM.Cells.Item(94, 11).Value = 639.2
$ws_BSM.Cells.Item(94, 12).Value = 1898.375
$ws_BSM.Cells.Item(94, 13).Value = -188.2
$ws_BSM.Cells.Item(94, 14).Value = -2800.375

# CRP row 31
$ws_CRP.Cells.Item(31, 8).Value = 1875.6757
$ws_CRP.Cells.Item(31, 9).Value = 1576.3334
$ws_CRP.Cells.Item(31, 11).Value = 1576.3334
$ws_CRP.Cells.Item(31, 13).Value = -1281.3334

# CRP row 34
$ws_CRP.Cells.Item(34, 8).Value = 1875.6757
$ws_CRP.Cells.Item(34, 9).Value = 1576.3334
$ws_CRP.Cells.Item(34, 11).Value = 1576.3334
$ws_CRP.Cells.Item(34, 13).Value = -1374.3334

# CRP row 39
$ws_CRP.Cells.Item(39, 8).Value = 29767.334
$ws_CRP.Cells.Item(39, 9).Value = 29767.334
$ws_CRP.Cells.Item(39, 11).Value = 29767.334
$ws_CRP.Cells.Item(39, 13).Value = -29376.334

# CRP row 49
$ws_CRP.Cells.Item(49, 8).Value = 29767.334
$ws_CRP.Cells.Item(49, 9).Value = 29767.334
$ws_CRP.Cells.Item(49, 11).Value = 29767.334
$ws_CRP.Cells.Item(49, 13).Value = -29585.334

# CRP row 62
$ws_CRP.Cells.Item(62, 8).Value = 127726.25
$ws_CRP.Cells.Item(62, 9).Value = 251502.5
$ws_CRP.Cells.Item(62, 10).Value = 3950
$ws_CRP.Cells.Item(62, 11).Value = 251502.5
$ws_CRP.Cells.Item(62, 12).Value = 3950
$ws_CRP.Cells.Item(62, 13).Value = -250878.5
$ws_CRP.Cells.Item(62, 14).Value = -5198

# CRP row 65
$ws_CRP.Cells.Item(65, 8).Value = 127726.25
$ws_CRP.Cells.Item(65, 9).Value = 251502.5
$ws_CRP.Cells.Item(65, 10).Value = 3950
$ws_CRP.Cells.Item(65, 11).Value = 1257512.5
$ws_CRP.Cells.Item(65, 12).Value = 19750
$ws_CRP.Cells.Item(65, 13).Value = -1254392.5
$ws_CRP.Cells.Item(65, 14).Value = -25990

# CRP row 132
$ws_CRP.Cells.Item(132, 8).Value = 331080.06
$ws_CRP.Cells.Item(132, 9).Value = 387222.4
$ws_CRP.Cells.Item(132, 10).Value = 3583
$ws_CRP.Cells.Item(132, 11).Value = 1161667.2
$ws_CRP.Cells.Item(132, 12).Value = 10749
$ws_CRP.Cells.Item(132, 13).Value = -1159137.2
$ws_CRP.Cells.Item(132, 14).Value = -15809

# CRP row 133
$ws_CRP.Cells.Item(133, 8).Value = 30326
$ws_CRP.Cells.Item(133, 10).Value = 30326
$ws_CRP.Cells.Item(133, 12).Value = 30326
$ws_CRP.Cells.Item(133, 14).Value = -35386

# CUL row 132
$ws_CUL.Cells.Item(132, 8).Value = 1711.1578
$ws_CUL.Cells.Item(132, 9).Value = 1140.8
$ws_CUL.Cells.Item(132, 10).Value = 1914.8572
$ws_CUL.Cells.Item(132, 11).Value = 10267.2
$ws_CUL.Cells.Item(132, 12).Value = 17233.7148
$ws_CUL.Cells.Item(132, 13).Value = -7737.199999999999
$ws_CUL.Cells.Item(132, 14).Value = -22293.7148

# GSM row 62
$ws_GSM.Cells.Item(62, 8).Value = 0
$ws_GSM.Cells.Item(62, 9).Value = 0
$ws_GSM.Cells.Item(62, 10).Value = 0
$ws_GSM.Cells.Item(62, 11).Value = 0
$ws_GSM.Cells.Item(62, 12).Value = 0
$ws_GSM.Cells.Item(62, 13).ClearContents()
$ws_GSM.Cells.Item(62, 14).ClearContents()

# GSM row 65
$ws_GSM.Cells.Item(65, 8).Value = 0
$ws_GSM.Cells.Item(65, 9).Value = 0
$ws_GSM.Cells.Item(65, 10).Value = 0
$ws_GSM.Cells.Item(65, 11).Value = 0
$ws_GSM.Cells.Item(65, 12).Value = 0
$ws_GSM.Cells.Item(65, 13).ClearContents()
$ws_GSM.Cells.Item(65, 14).ClearContents()

# LTW row 7
$ws_LTW.Cells.Item(7, 8).Value = 3977.3333
$ws_LTW.Cells.Item(7, 9).Value = 4028.5715
$ws_LTW.Cells.Item(7, 10).Value = 3932.5
$ws_LTW.Cells.Item(7, 11).Value = 4028.5715
$ws_LTW.Cells.Item(7, 12).Value = 3932.5
$ws_LTW.Cells.Item(7, 13).Value = -3916.5715
$ws_LTW.Cells.Item(7, 14).Value = -4156.5

# LTW row 22
$ws_LTW.Cells.Item(22, 8).Value = 508.45456
$ws_LTW.Cells.Item(22, 9).Value = 460.14285
$ws_LTW.Cells.Item(22, 10).Value = 593
$ws_LTW.Cells.Item(22, 11).Value = 460.14285
$ws_LTW.Cells.Item(22, 12).Value = 593
$ws_LTW.Cells.Item(22, 13).Value = -165.14285
$ws_LTW.Cells.Item(22, 14).Value = -1183

# LTW row 27
$ws_LTW.Cells.Item(27, 8).Value = 508.45456
$ws_LTW.Cells.Item(27, 9).Value = 460.14285
$ws_LTW.Cells.Item(27, 10).Value = 593
$ws_LTW.Cells.Item(27, 11).Value = 460.14285
$ws_LTW.Cells.Item(27, 12).Value = 593
$ws_LTW.Cells.Item(27, 13).Value = -353.14285
$ws_LTW.Cells.Item(27, 14).Value = -807

# LTW row 40
$ws_LTW.Cells.Item(40, 8).Value = 6279.6665
$ws_LTW.Cells.Item(40, 9).Value = 6335.6
$ws_LTW.Cells.Item(40, 10).Value = 6000
$ws_LTW.Cells.Item(40, 11).Value = 6335.6
$ws_LTW.Cells.Item(40, 12).Value = 6000
$ws_LTW.Cells.Item(40, 13).Value = -6199.6
$ws_LTW.Cells.Item(40, 14).Value = -6272

# LTW row 64
$ws_LTW.Cells.Item(64, 8).Value = 21000
$ws_LTW.Cells.Item(64, 9).Value = 12000
$ws_LTW.Cells.Item(64, 11).Value = 12000
$ws_LTW.Cells.Item(64, 13).Value = -11775

# LTW row 67
$ws_LTW.Cells.Item(67, 8).Value = 21000
$ws_LTW.Cells.Item(67, 9).Value = 12000
$ws_LTW.Cells.Item(67, 11).Value = 12000
$ws_LTW.Cells.Item(67, 13).Value = -11220

# LTW row 68
$ws_LTW.Cells.Item(68, 8).Value = 3285.7144
$ws_LTW.Cells.Item(68, 9).Value = 3750
$ws_LTW.Cells.Item(68, 10).Value = 2666.6667
$ws_LTW.Cells.Item(68, 11).Value = 3750
$ws_LTW.Cells.Item(68, 12).Value = 2666.6667
$ws_LTW.Cells.Item(68, 13).Value = -3001
$ws_LTW.Cells.Item(68, 14).Value = -4164.6667

# LTW row 71
$ws_LTW.Cells.Item(71, 8).Value = 3285.7144
$ws_LTW.Cells.Item(71, 9).Value = 3750
$ws_LTW.Cells.Item(71, 10).Value = 2666.6667
$ws_LTW.Cells.Item(71, 11).Value = 18750
$ws_LTW.Cells.Item(71, 12).Value = 13333.3335
$ws_LTW.Cells.Item(71, 13).Value = -15006
$ws_LTW.Cells.Item(71, 14).Value = -20821.3335

# LTW row 76
$ws_LTW.Cells.Item(76, 8).Value = 7199.75
$ws_LTW.Cells.Item(76, 10).Value = 7799.7144
$ws_LTW.Cells.Item(76, 12).Value = 7799.7144
$ws_LTW.Cells.Item(76, 14).Value = -8475.714400000001

# LTW row 79
$ws_LTW.Cells.Item(79, 8).Value = 7199.75
$ws_LTW.Cells.Item(79, 10).Value = 7799.7144
$ws_LTW.Cells.Item(79, 12).Value = 7799.7144
$ws_LTW.Cells.Item(79, 14).Value = -10139.7144

# LTW row 97
$ws_LTW.Cells.Item(97, 8).Value = 22341.277
$ws_LTW.Cells.Item(97, 10).Value = 22341.277
$ws_LTW.Cells.Item(97, 12).Value = 22341.277
$ws_LTW.Cells.Item(97, 14).Value = -24323.277

# LTW row 112
$ws_LTW.Cells.Item(112, 8).Value = 27500
$ws_LTW.Cells.Item(112, 10).Value = 27500
$ws_LTW.Cells.Item(112, 12).Value = 27500
$ws_LTW.Cells.Item(112, 14).Value = -30454

# LTW row 126
$ws_LTW.Cells.Item(126, 8).Value = 3977.3333
$ws_LTW.Cells.Item(126, 9).Value = 4028.5715
$ws_LTW.Cells.Item(126, 10).Value = 3932.5
$ws_LTW.Cells.Item(126, 11).Value = 12085.7145
$ws_LTW.Cells.Item(126, 12).Value = 11797.5
$ws_LTW.Cells.Item(126, 13).Value = -9615.7145
$ws_LTW.Cells.Item(126, 14).Value = -16737.5

# LTW row 140
$ws_LTW.Cells.Item(140, 8).Value = 39356.445
$ws_LTW.Cells.Item(140, 10).Value = 39356.445
$ws_LTW.Cells.Item(140, 12).Value = 39356.445
$ws_LTW.Cells.Item(140, 14).Value = -49716.445

# WVR row 62
$ws_WVR.Cells.Item(62, 8).Value = 4600

# WVR row 65
$ws_WVR.Cells.Item(65, 8).Value = 4600

# WVR row 96
$ws_WVR.Cells.Item(96, 8).Value = 1466.6666
$ws_WVR.Cells.Item(96, 9).Value = 0
$ws_WVR.Cells.Item(96, 10).Value = 1466.6666
$ws_WVR.Cells.Item(96, 11).Value = 0
$ws_WVR.Cells.Item(96, 12).Value = 1466.6666
$ws_WVR.Cells.Item(96, 13).ClearContents()
$ws_WVR.Cells.Item(96, 14).Value = -4212.6666
